# ExamCenterDetails.xlsx — "Latest Code & Handled Stale Element reference
# Exception": the STAGE sheet's sample login/exam row (row 2) now points at
# a fresh exam-center Location / ExamName / ScheduleName, and the three
# "last run" schedule ids below it (E3:E5) were bumped to new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STAGE")

# --- Row 2: Location / ExamName / ScheduleName -----------------------
# None of these look like numbers, so a plain .Value assignment keeps
# them as text (t="s") and leaves the existing cell formatting alone.
$ws.Range("A2").Value = "ECLocation50954"
$ws.Range("H2").Value = "FPK12Exam53770"
$ws.Range("I2").Value = "FPK12Schedule28586"

# --- E3:E5: numeric-looking id strings --------------------------------
# These values (e.g. "45044") must stay text cells, matching the
# original data. Assigning a digit-only string straight to .Value would
# get auto-coerced to a number by Excel, so build each one as a text
# formula result in a scratch cell, then bring over just the value via
# PasteSpecial (values only) — that preserves the text type without
# disturbing the destination cell's existing style. The scratch column
# is removed afterwards so it leaves no trace.
$scratch = $ws.Range("ZZ1")

$scratch.Formula = '="45044"'
$scratch.Copy()
$ws.Range("E3").PasteSpecial(-4163)

$scratch.Formula = '="59184"'
$scratch.Copy()
$ws.Range("E4").PasteSpecial(-4163)

$scratch.Formula = '="42692"'
$scratch.Copy()
$ws.Range("E5").PasteSpecial(-4163)

$scratch.EntireColumn.Delete()
